$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptos list values.
# Column D ("Price") values are forced to remain plain text (they look numeric,
# e.g. "3.00" / "69.986.88", and Excel would otherwise silently convert them to
# floating point numbers and normalize/round their textual representation).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.986.88"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.53%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.546.63"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.71%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "587.31"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.04%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "185.41"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.37%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.535.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.89%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.615"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.37%  "

$ws.Range("E9").Value = "  -0.02%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.198"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.68%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.645"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.05%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "54.27"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.11%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000306"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.31%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.49"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.25%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.106.68"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.98%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "19.38"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.21%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "69.978.26"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.60%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.519.32"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.82%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.45"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.02%  "

$ws.Range("E20").Value = "  -1.26%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "539.14"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +9.20%  "

$ws.Range("E22").Value = "  -2.90%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "17.95"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -8.48%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.60"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.93%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.87"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.24%  "

$ws.Range("E26").Value = "  -0.55%  "

$ws.Range("B27").Value = "ImmutableX"
$ws.Range("C27").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.15%  "

$ws.Range("B28").Value = "RenderToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "11.20"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.84%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.13"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.29%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "32.13"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.35%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.29"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.33%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "12.47"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.04%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "64.51"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.94%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.113"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.09%  "

$ws.Range("B35").Value = "Fetch.AI"
$ws.Range("C35").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.25"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +8.18%  "

$ws.Range("B36").Value = "Bittensor"
$ws.Range("C36").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "548.04"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.61%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.414"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.07%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "38.34"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.43%  "

$ws.Range("E39").Value = "  -0.11%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0₃0764"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -6.31%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.135"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.08%  "

$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.40"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.17%  "

$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.345.80"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.30%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.10"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -8.77%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.58"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.92%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.97"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.18%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0441"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.57%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.16"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -6.74%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.136"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.74%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.998"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.14%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "137.21"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.36%  "
